$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (the first data row) in place before removing the rest.
$ws.Range("C2").Value = "May 3"
$ws.Range("D2").Value = "12:00 AM"
# Leading apostrophe keeps "10" stored as text (matches the rest of the
# sheet, which is entirely text-typed) instead of Excel auto-coercing it
# to a number.
$ws.Range("E2").Value = "'10"

# Remove the remaining data rows (old rows 3-8), shrinking the sheet down
# to just the header row and the single remaining data row.
$ws.Rows("3:8").Delete() | Out-Null
